$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style/format as the existing header cell (H1) to the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF), rows 2-20
$data = @{
    2  = @(6, 6)
    3  = @(8, 8)
    4  = @(8, 9)
    5  = @(8, 8)
    6  = @(9, 9)
    7  = @(7, 7)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(6, 6)
    11 = @(6, 7)
    12 = @(8, 8)
    13 = @(6, 7)
    14 = @(7, 7)
    15 = @(5, 5)
    16 = @(7, 7)
    17 = @(6, 7)
    18 = @(7, 7)
    19 = @(5, 5)
    20 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
